# Updates the "Forecast Comparison" sheet (new forecast run: every week's
# date rolls forward by one, W1's old row is dropped, and a new W16 week is
# appended) and the "Summary" sheet's derived metrics, matching the
# "Fixed update to excel issue" forecast refresh.

$wb = $excel.ActiveWorkbook
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison: rows 2-17 (B=Week_Start_Date, D..H=forecast cols) ---
# Each row keeps its Week label (A) / ASIN (C) / Product Title (I) /
# is_holiday_week (J) as-is; only the date + numeric forecast columns move.

$forecastRows = @(
    @{ Row = 2;  B = "2025-02-02"; D = 17; E = 42; F = 50; G = 57; H = 69 },
    @{ Row = 3;  B = "2025-02-09"; D = 17; E = 43; F = 52; G = 60; H = 73 },
    @{ Row = 4;  B = "2025-02-16"; D = 17; E = 43; F = 51; G = 60; H = 73 },
    @{ Row = 5;  B = "2025-02-23"; D = 16; E = 43; F = 51; G = 59; H = 73 },
    @{ Row = 6;  B = "2025-03-02"; D = 15; E = 43; F = 52; G = 61; H = 76 },
    @{ Row = 7;  B = "2025-03-09"; D = 14; E = 41; F = 50; G = 58; H = 72 },
    @{ Row = 8;  B = "2025-03-16"; D = 15; E = 43; F = 52; G = 62; H = 78 },
    @{ Row = 9;  B = "2025-03-23"; D = 16; E = 45; F = 54; G = 66; H = 85 },
    @{ Row = 10; B = "2025-03-30"; D = 16; E = 41; F = 50; G = 59; H = 74 },
    @{ Row = 11; B = "2025-04-06"; D = 16; E = 41; F = 50; G = 62; H = 79 },
    @{ Row = 12; B = "2025-04-13"; D = 16; E = 42; F = 52; G = 63; H = 81 },
    @{ Row = 13; B = "2025-04-20"; D = 17; E = 42; F = 51; G = 63; H = 80 },
    @{ Row = 14; B = "2025-04-27"; D = 17; E = 42; F = 51; G = 62; H = 79 },
    @{ Row = 15; B = "2025-05-04"; D = 16; E = 39; F = 48; G = 59; H = 76 },
    @{ Row = 16; B = "2025-05-11"; D = 16; E = 39; F = 48; G = 59; H = 77 },
    @{ Row = 17; B = "2025-05-18"; D = 16; E = 39; F = 47; G = 59; H = 77 }
)

foreach ($r in $forecastRows) {
    $row = $r.Row

    # Week_Start_Date is stored as plain text "yyyy-MM-dd", not a date
    # serial -- force text so Excel doesn't auto-convert it to a date.
    $cellB = $wsForecast.Cells.Item($row, 2)
    $cellB.NumberFormat = "@"
    $cellB.Value = $r.B

    $wsForecast.Cells.Item($row, 4).Value = $r.D
    $wsForecast.Cells.Item($row, 5).Value = $r.E
    $wsForecast.Cells.Item($row, 6).Value = $r.F
    $wsForecast.Cells.Item($row, 7).Value = $r.G
    $wsForecast.Cells.Item($row, 8).Value = $r.H
}

# --- Summary: column B values are all stored as text, including the
# purely-numeric-looking ones, so force text on every updated cell. ---

$summaryRows = @(
    @{ Row = 2;  Value = "2022-12-25 to 2025-01-26" },
    @{ Row = 4;  Value = "240" },
    @{ Row = 5;  Value = "75" },
    @{ Row = 6;  Value = "72" },
    @{ Row = 7;  Value = "62" },
    @{ Row = 8;  Value = "8181 units" },
    @{ Row = 9;  Value = "253" },
    @{ Row = 10; Value = "125" },
    @{ Row = 11; Value = "66" },
    @{ Row = 12; Value = "17" },
    @{ Row = 13; Value = "2025-02-09" },
    @{ Row = 14; Value = "14" },
    @{ Row = 15; Value = "2025-03-09" }
)

foreach ($r in $summaryRows) {
    $cell = $wsSummary.Cells.Item($r.Row, 2)
    $cell.NumberFormat = "@"
    $cell.Value = $r.Value
}
